# Fruta / hortaliza, semanal
# Insert a new weekly price record at row 122 for
# "Femacal de La Calera" - Arándano (blue). This pushes the existing
# rows 122:143 down to 123:144 (dimension grows from A1:T143 to A1:T144).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 122:143 down by inserting a fresh row at 122.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row with this week's record.
$ws.Range("A122").Value = 3
$ws.Range("B122").Value = "Femacal de La Calera"
$ws.Range("C122").Value = "Coquimbo"
$ws.Range("D122").Value = 44505
$ws.Range("E122").Value = 5
$ws.Range("F122").Value = "Fruta"
$ws.Range("G122").Value = 100101
$ws.Range("H122").Value = "Berries"
$ws.Range("I122").Value = 100101001
$ws.Range("J122").Value = "Arándano (blue)"
$ws.Range("K122").Value = "Sin especificar"
$ws.Range("L122").Value = "Primera"
$ws.Range("M122").Value = 38
$ws.Range("N122").Value = 10000
$ws.Range("O122").Value = 10000
$ws.Range("P122").Value = 10000
$ws.Range("Q122").Value = "$/bandeja 2 kilos"
$ws.Range("R122").Value = "Provincia de Quillota"
$ws.Range("S122").Value = 5000
$ws.Range("T122").Value = 2
